$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "ISSUES" section - replace the "Couldn't resolve..." paragraph
# with the "Resolved - 5/15/2023..." paragraph, and drop the blank paragraph
# that used to precede it.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("Couldn" + [char]0x2019 + "t resolve.  ", $false, $false, $false, $false, $false, $true, 1, $false, "Resolved " + [char]0x2013 + " 5/15/2023: ", 2)
$d.Content.Find.Execute(" remains NaN when passing Population to it.", $false, $false, $false, $false, $false, $true, 1, $false, " domain is expecting an array with 2 elements; ", 2)

$i = 0
$target = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -match "^Resolved") {
        $target = $i
    }
}
$blank = $d.Paragraphs.Item($target - 1)
$blank.Range.Delete()

# ---------------------------------------------------------------------------
# Change 2: "THINGS TO DO" section - mark the yScale item as done, and add
# two new list items after it.
# ---------------------------------------------------------------------------

$i = 0
$target = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -match "^Find and fix the issue") {
        $target = $i
    }
}

$para = $d.Paragraphs.Item($target)
$paraRange = $para.Range
$dot = $d.Range($paraRange.End - 2, $paraRange.End - 1)
$dot.Text = " - done 5/15/2023"

$nextPara = $d.Paragraphs.Item($target + 1)
$insPoint = $nextPara.Range.Start
$ins = $d.Range($insPoint, $insPoint)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="202122"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F8F9FA"/></w:rPr>' + `
  '<w:t>Fix issue with all state populations not appearing on bar chart when no state is selected (setting this as a default</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Fix color legend so that colors other than ORANGE are used (when different states are selected)</w:t></w:r></w:p>' + `
  '<w:p/>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$ins.InsertXML($xml)

# Remove the stray blank paragraph left behind by InsertXML's trailing <w:p/>
# (it sits right before the paragraph that used to follow our insertion point).
$i = 0
$target = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -match "^Fix color legend") {
        $target = $i
    }
}
$stray = $d.Paragraphs.Item($target + 1)
$stray.Range.Delete()

Write-Output "done"
